# TestData update by shany Mohan Dhas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (TC_OpenMrs_Delete_a_Patient) - FirstName
$ws.Range("D5").Value = "Zig zag"

# Row 2 (TC_OpenMrs_Register_a_Patient) - FirstName, MiddleName, FamilyName
$ws.Range("D2").Value = "Letin"
$ws.Range("E2").Value = "John"
$ws.Range("F2").Value = "K"

# Leave the active selection on F2, matching the saved workbook view
$ws.Range("F2").Select()
